$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 66670660
$ws.Range("J76").Value = 5358.143
$ws.Range("L76").Value = 5358.143
$ws.Range("N76").Value = -5988.143

$ws.Range("H79").Value = 66670660
$ws.Range("J79").Value = 5358.143
$ws.Range("L79").Value = 5358.143
$ws.Range("N79").Value = -7542.143

$ws.Range("H98").Value = 1370.7
$ws.Range("I98").Value = 1267.1111
$ws.Range("J98").Value = 2303
$ws.Range("K98").Value = 1267.1111
$ws.Range("L98").Value = 2303
$ws.Range("M98").Value = 230.8888999999999
$ws.Range("N98").Value = -5299

$ws.Range("H100").Value = 3863.8948
$ws.Range("I100").Value = 1934
$ws.Range("J100").Value = 7172.2856
$ws.Range("K100").Value = 1934
$ws.Range("L100").Value = 7172.2856
$ws.Range("M100").Value = -1393
$ws.Range("N100").Value = -8254.285599999999

$ws.Range("H122").Value = 1370.7
$ws.Range("I122").Value = 1267.1111
$ws.Range("J122").Value = 2303
$ws.Range("K122").Value = 3801.3333
$ws.Range("L122").Value = 6909
$ws.Range("M122").Value = -1351.3333
$ws.Range("N122").Value = -11809

$ws.Range("H129").Value = 2570.5789
$ws.Range("I129").Value = 386.2143
$ws.Range("J129").Value = 8686.799999999999
$ws.Range("K129").Value = 1158.6429
$ws.Range("L129").Value = 26060.4
$ws.Range("M129").Value = 3841.3571
$ws.Range("N129").Value = -36060.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 41476.84
$ws.Range("I74").Value = 72989.17999999999
$ws.Range("J74").Value = 1370.2273
$ws.Range("K74").Value = 72989.17999999999
$ws.Range("L74").Value = 1370.2273
$ws.Range("M74").Value = -72115.17999999999
$ws.Range("N74").Value = -3118.2273

$ws.Range("H77").Value = 41476.84
$ws.Range("I77").Value = 72989.17999999999
$ws.Range("J77").Value = 1370.2273
$ws.Range("K77").Value = 364945.9
$ws.Range("L77").Value = 6851.136500000001
$ws.Range("M77").Value = -360577.9
$ws.Range("N77").Value = -15587.1365

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1346
$ws.Range("I99").Value = 1109.7142
$ws.Range("K99").Value = 1109.7142
$ws.Range("M99").Value = 388.2858000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32281.47
$ws.Range("I31").Value = 63521.75
$ws.Range("J31").Value = 4512.3335
$ws.Range("K31").Value = 63521.75
$ws.Range("L31").Value = 4512.3335
$ws.Range("M31").Value = -63226.75
$ws.Range("N31").Value = -5102.3335

$ws.Range("H34").Value = 32281.47
$ws.Range("I34").Value = 63521.75
$ws.Range("J34").Value = 4512.3335
$ws.Range("K34").Value = 63521.75
$ws.Range("L34").Value = 4512.3335
$ws.Range("M34").Value = -63319.75
$ws.Range("N34").Value = -4916.3335

$ws.Range("H37").Value = 70028.5
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 70028.5
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 70028.5
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -70242.5

$ws.Range("H41").Value = 10377.883
$ws.Range("I41").Value = 4619.6665
$ws.Range("J41").Value = 11611.786
$ws.Range("K41").Value = 4619.6665
$ws.Range("L41").Value = 11611.786
$ws.Range("M41").Value = -4191.6665
$ws.Range("N41").Value = -12467.786

$ws.Range("H141").Value = 49109.89
$ws.Range("J141").Value = 49109.89
$ws.Range("L141").Value = 49109.89
$ws.Range("N141").Value = -59469.89

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1801.0625
$ws.Range("I46").Value = 183.33333
$ws.Range("J46").Value = 2174.3845
$ws.Range("K46").Value = 549.99999
$ws.Range("L46").Value = 6523.1535
$ws.Range("M46").Value = -458.99999
$ws.Range("N46").Value = -6705.1535

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 20086000
$ws.Range("I11").Value = 40167000
$ws.Range("J11").Value = 4999
$ws.Range("K11").Value = 40167000
$ws.Range("L11").Value = 4999
$ws.Range("M11").Value = -40166861
$ws.Range("N11").Value = -5277

$ws.Range("H21").Value = 50002.332
$ws.Range("I21").Value = 10000
$ws.Range("J21").Value = 58002.8
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 58002.8
$ws.Range("M21").Value = -9827
$ws.Range("N21").Value = -58348.8

$ws.Range("H30").Value = 50002.332
$ws.Range("I30").Value = 10000
$ws.Range("J30").Value = 58002.8
$ws.Range("K30").Value = 10000
$ws.Range("L30").Value = 58002.8
$ws.Range("M30").Value = -9895
$ws.Range("N30").Value = -58212.8

$ws.Range("H64").Value = 20300
$ws.Range("I64").Value = 10000
$ws.Range("J64").Value = 22875
$ws.Range("K64").Value = 10000
$ws.Range("L64").Value = 22875
$ws.Range("M64").Value = -9752
$ws.Range("N64").Value = -23371

$ws.Range("H67").Value = 20300
$ws.Range("I67").Value = 10000
$ws.Range("J67").Value = 22875
$ws.Range("K67").Value = 10000
$ws.Range("L67").Value = 22875
$ws.Range("M67").Value = -9142
$ws.Range("N67").Value = -24591

$ws.Range("H70").Value = 4412.125
$ws.Range("I70").Value = 4311.56
$ws.Range("K70").Value = 4311.56
$ws.Range("M70").Value = -4041.56

$ws.Range("H73").Value = 4412.125
$ws.Range("I73").Value = 4311.56
$ws.Range("K73").Value = 4311.56
$ws.Range("M73").Value = -3375.56

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 18900
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 18900
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 18900
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -19776

$ws.Range("H68").Value = 2358
$ws.Range("I68").Value = 1788.7778
$ws.Range("J68").Value = 2870.3
$ws.Range("K68").Value = 1788.7778
$ws.Range("L68").Value = 2870.3
$ws.Range("M68").Value = -1039.7778
$ws.Range("N68").Value = -4368.3

$ws.Range("H71").Value = 2358
$ws.Range("I71").Value = 1788.7778
$ws.Range("J71").Value = 2870.3
$ws.Range("K71").Value = 8943.889000000001
$ws.Range("L71").Value = 14351.5
$ws.Range("M71").Value = -5199.889000000001
$ws.Range("N71").Value = -21839.5

$ws.Range("H136").Value = 295347.72
$ws.Range("I136").Value = 417249.28
$ws.Range("J136").Value = 2783.9
$ws.Range("K136").Value = 1251747.84
$ws.Range("L136").Value = 8351.700000000001
$ws.Range("M136").Value = -1249197.84
$ws.Range("N136").Value = -13451.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1600
$ws.Range("I126").Value = 1250
$ws.Range("J126").Value = 1950
$ws.Range("K126").Value = 3750
$ws.Range("L126").Value = 5850
$ws.Range("M126").Value = -1280
$ws.Range("N126").Value = -10790
